$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2=44420, J2=700, K2=27000, L2=29000, M2=28000, P2=1120
$ws.Range("D2").Value = 44420
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 27000
$ws.Range("L2").Value = 29000
$ws.Range("M2").Value = 28000
$ws.Range("P2").Value = 1120

# Row 3: D3=44419, J3=600, K3=27000, L3=29000, M3=28000, P3=1120
$ws.Range("D3").Value = 44419
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 27000
$ws.Range("L3").Value = 29000
$ws.Range("M3").Value = 28000
$ws.Range("P3").Value = 1120

# Row 4: D4=44377, J4=500, K4=26000, L4=28000, M4=27000, P4=1080
$ws.Range("D4").Value = 44377
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 26000
$ws.Range("L4").Value = 28000
$ws.Range("M4").Value = 27000
$ws.Range("P4").Value = 1080

# Row 5: D5=44357, J5=340, K5=28000, L5=30000, M5=29000, P5=1160
$ws.Range("D5").Value = 44357
$ws.Range("J5").Value = 340
$ws.Range("K5").Value = 28000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 29000
$ws.Range("P5").Value = 1160

# Row 6: D6=44384, J6=400, K6=26000, L6=28000, M6=27000, P6=1080
$ws.Range("D6").Value = 44384
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 26000
$ws.Range("L6").Value = 28000
$ws.Range("M6").Value = 27000
$ws.Range("P6").Value = 1080

# Row 7: D7=44406, J7=600, K7=26000, L7=28000, M7=27000, P7=1080
$ws.Range("D7").Value = 44406
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 26000
$ws.Range("L7").Value = 28000
$ws.Range("M7").Value = 27000
$ws.Range("P7").Value = 1080

# Row 8: D8=44363, J8=240, K8=28000, L8=30000, M8=29000, P8=1160
$ws.Range("D8").Value = 44363
$ws.Range("J8").Value = 240
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 29000
$ws.Range("P8").Value = 1160

# Row 9: D9=44343, J9=200, K9=26000, L9=28000, M9=27000, P9=1080
$ws.Range("D9").Value = 44343
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 26000
$ws.Range("L9").Value = 28000
$ws.Range("M9").Value = 27000
$ws.Range("P9").Value = 1080

# Row 10: D10=44490, J10=500, K10=16000, L10=18000, M10=17000, P10=680
$ws.Range("D10").Value = 44490
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17000
$ws.Range("P10").Value = 680

# Row 11: D11=44503, J11=400, K11=11000, L11=13000, M11=12000, P11=480
$ws.Range("D11").Value = 44503
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 12000
$ws.Range("P11").Value = 480

# Row 12: D12=44448, J12=400, K12=28000, L12=30000, M12=29000, P12=1160
$ws.Range("D12").Value = 44448
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 28000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29000
$ws.Range("P12").Value = 1160

# Row 13: D13=44497, J13=500, K13=13000, L13=15000, M13=14000, P13=560
$ws.Range("D13").Value = 44497
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14000
$ws.Range("P13").Value = 560

# Row 14: D14=44427, J14=300, K14=28000, L14=30000, M14=29000, P14=1160
$ws.Range("D14").Value = 44427
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 28000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29000
$ws.Range("P14").Value = 1160

# Row 15: D15=44482, J15=500, K15=18000, L15=20000, M15=19000, P15=760
$ws.Range("D15").Value = 44482
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19000
$ws.Range("P15").Value = 760

# Row 16: D16=44412, J16=600, K16=25000, L16=27000, M16=26000, P16=1040
$ws.Range("D16").Value = 44412
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 27000
$ws.Range("M16").Value = 26000
$ws.Range("P16").Value = 1040

# Row 17: D17=44483, J17=300, K17=18000, L17=20000, M17=19000, P17=760
$ws.Range("D17").Value = 44483
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19000
$ws.Range("P17").Value = 760

# Row 18: D18=44370, J18=400, K18=27000, L18=28000, M18=27500, P18=1100
$ws.Range("D18").Value = 44370
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 27000
$ws.Range("L18").Value = 28000
$ws.Range("M18").Value = 27500
$ws.Range("P18").Value = 1100

# Row 19: D19=44469, J19=600, K19=22000, L19=24000, M19=23000, P19=920
$ws.Range("D19").Value = 44469
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 22000
$ws.Range("L19").Value = 24000
$ws.Range("M19").Value = 23000
$ws.Range("P19").Value = 920

# Row 20: D20=44349, J20=600, K20=26000, L20=28000, M20=27000, P20=1080
$ws.Range("D20").Value = 44349
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 26000
$ws.Range("L20").Value = 28000
$ws.Range("M20").Value = 27000
$ws.Range("P20").Value = 1080

# Row 21: D21=44398, J21=500, K21=26000, L21=28000, M21=27000, P21=1080
$ws.Range("D21").Value = 44398
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 26000
$ws.Range("L21").Value = 28000
$ws.Range("M21").Value = 27000
$ws.Range("P21").Value = 1080

# Row 22: D22=44435, J22=900, K22=28000, L22=30000, M22=29000, P22=1160
$ws.Range("D22").Value = 44435
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 28000
$ws.Range("L22").Value = 30000
$ws.Range("M22").Value = 29000
$ws.Range("P22").Value = 1160

# Row 23: D23=44350, J23=700, K23=28000, L23=30000, M23=29000, P23=1160
$ws.Range("D23").Value = 44350
$ws.Range("J23").Value = 700
$ws.Range("K23").Value = 28000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 29000
$ws.Range("P23").Value = 1160

# Row 24: D24=44413, J24=700, K24=26000, L24=28000, M24=27000, P24=1080
$ws.Range("D24").Value = 44413
$ws.Range("J24").Value = 700
$ws.Range("K24").Value = 26000
$ws.Range("L24").Value = 28000
$ws.Range("M24").Value = 27000
$ws.Range("P24").Value = 1080

# Row 25: D25=44405, J25=500, K25=26000, L25=28000, M25=27000, P25=1080
$ws.Range("D25").Value = 44405
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 26000
$ws.Range("L25").Value = 28000
$ws.Range("M25").Value = 27000
$ws.Range("P25").Value = 1080

# Row 26: D26=44434, J26=500, K26=28000, L26=30000, M26=29000, P26=1160
$ws.Range("D26").Value = 44434
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 28000
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = 29000
$ws.Range("P26").Value = 1160

# Row 27: D27=44476, J27=500, K27=23000, L27=24000, M27=23500, P27=940
$ws.Range("D27").Value = 44476
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 23000
$ws.Range("L27").Value = 24000
$ws.Range("M27").Value = 23500
$ws.Range("P27").Value = 940

# Row 28: D28=44433, J28=400, K28=28000, L28=30000, M28=29000, P28=1160
$ws.Range("D28").Value = 44433
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 29000
$ws.Range("P28").Value = 1160

# Row 29: D29=44356, J29=300, K29=26000, L29=28000, M29=27000, P29=1080
$ws.Range("D29").Value = 44356
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 26000
$ws.Range("L29").Value = 28000
$ws.Range("M29").Value = 27000
$ws.Range("P29").Value = 1080

# Row 30: D30=44455, J30=800, K30=28000, L30=30000, M30=29000, P30=1160
$ws.Range("D30").Value = 44455
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 28000
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = 29000
$ws.Range("P30").Value = 1160

# Row 31: D31=44461, J31=500, K31=23000, L31=25000, M31=24000, P31=960
$ws.Range("D31").Value = 44461
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 23000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 24000
$ws.Range("P31").Value = 960

# Row 32: D32=44462, J32=400, K32=22000, L32=23000, M32=22500, P32=900
$ws.Range("D32").Value = 44462
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 22000
$ws.Range("L32").Value = 23000
$ws.Range("M32").Value = 22500
$ws.Range("P32").Value = 900

# Row 33: D33=44364, J33=200, K33=28000, L33=30000, M33=29000, P33=1160
$ws.Range("D33").Value = 44364
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 28000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = 29000
$ws.Range("P33").Value = 1160

# Row 34: D34=44426, J34=400, K34=28000, L34=30000, M34=29000, P34=1160
$ws.Range("D34").Value = 44426
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 28000
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = 29000
$ws.Range("P34").Value = 1160

# Row 35: D35=44504, J35=600, K35=11000, L35=13000, M35=12000, P35=480
$ws.Range("D35").Value = 44504
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 12000
$ws.Range("P35").Value = 480

# Row 36: D36=44447, J36=600, K36=28000, L36=30000, M36=29000, P36=1160
$ws.Range("D36").Value = 44447
$ws.Range("J36").Value = 600
$ws.Range("K36").Value = 28000
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = 29000
$ws.Range("P36").Value = 1160

# Row 37: D37=44391, J37=100, K37=26000, L37=28000, M37=27000, P37=1080
$ws.Range("D37").Value = 44391
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 26000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 27000
$ws.Range("P37").Value = 1080

# Row 38: D38=44454, J38=1000, K38=28000, L38=30000, M38=29000, P38=1160
$ws.Range("D38").Value = 44454
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 28000
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = 29000
$ws.Range("P38").Value = 1160

# Row 39: D39=44371, J39=500, K39=28000, L39=30000, M39=29000, P39=1160
$ws.Range("D39").Value = 44371
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 28000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 29000
$ws.Range("P39").Value = 1160

# Row 40: D40=44489, J40=400, K40=18000, L40=20000, M40=19000, P40=760
$ws.Range("D40").Value = 44489
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 18000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 19000
$ws.Range("P40").Value = 760

# Row 41: D41=44399, J41=400, K41=26000, L41=28000, M41=27000, P41=1080
$ws.Range("D41").Value = 44399
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 26000
$ws.Range("L41").Value = 28000
$ws.Range("M41").Value = 27000
$ws.Range("P41").Value = 1080

# Row 42: D42=44441, J42=700, K42=28000, L42=30000, M42=29000, P42=1160
$ws.Range("D42").Value = 44441
$ws.Range("J42").Value = 700
$ws.Range("K42").Value = 28000
$ws.Range("L42").Value = 30000
$ws.Range("M42").Value = 29000
$ws.Range("P42").Value = 1160

# Row 43: D43=44475, J43=1000, K43=22000, L43=24000, M43=23000, P43=920
$ws.Range("D43").Value = 44475
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 22000
$ws.Range("L43").Value = 24000
$ws.Range("M43").Value = 23000
$ws.Range("P43").Value = 920

# Row 44: D44=44468, J44=500, K44=23000, L44=25000, M44=24000, P44=960
$ws.Range("D44").Value = 44468
$ws.Range("J44").Value = 500
$ws.Range("K44").Value = 23000
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 24000
$ws.Range("P44").Value = 960

# Row 45: D45=44385, J45=500, K45=26000, L45=28000, M45=27000, P45=1080
$ws.Range("D45").Value = 44385
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 26000
$ws.Range("L45").Value = 28000
$ws.Range("M45").Value = 27000
$ws.Range("P45").Value = 1080

# Row 46: D46=44392, J46=100, K46=26000, L46=28000, M46=27000, P46=1080
$ws.Range("D46").Value = 44392
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 26000
$ws.Range("L46").Value = 28000
$ws.Range("M46").Value = 27000
$ws.Range("P46").Value = 1080
